$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.345.00'
$ws.Range('E2').Value = '  -2.53%  '
$ws.Range('D3').Value = '3.001.27'
$ws.Range('E3').Value = '  -2.20%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'538.20"
$ws.Range('E5').Value = '  +0.07%  '
$ws.Range('D6').Value = "'135.68"
$ws.Range('E6').Value = '  +1.70%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').Value = '2.995.24'
$ws.Range('E8').Value = '  -2.16%  '
$ws.Range('D9').Value = "'0.496"
$ws.Range('E9').Value = '  +0.38%  '
$ws.Range('E10').Value = '  -3.30%  '
$ws.Range('E11').Value = '  +0.67%  '
$ws.Range('E12').Value = '  -1.10%  '
$ws.Range('E13').Value = '  -0.99%  '
$ws.Range('D14').Value = "'34.25"
$ws.Range('E14').Value = '  -0.12%  '
$ws.Range('D15').Value = '3.492.92'
$ws.Range('E15').Value = '  -1.97%  '
$ws.Range('E16').Value = '  -0.52%  '
$ws.Range('D17').Value = '61.428.87'
$ws.Range('E17').Value = '  -2.34%  '
$ws.Range('D18').Value = '3.004.21'
$ws.Range('E18').Value = '  -2.04%  '
$ws.Range('D19').Value = "'6.63"
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('D20').Value = "'468.46"
$ws.Range('E20').Value = '  -2.93%  '
$ws.Range('D21').Value = "'13.27"
$ws.Range('E21').Value = '  -0.48%  '
$ws.Range('D22').Value = "'0.677"
$ws.Range('E22').Value = '  -2.71%  '
$ws.Range('D23').Value = "'6.97"
$ws.Range('E23').Value = '  -1.84%  '
$ws.Range('D24').Value = "'79.77"
$ws.Range('E24').Value = '  +0.94%  '
$ws.Range('D25').Value = "'12.05"
$ws.Range('E25').Value = '  -0.72%  '
$ws.Range('D26').Value = "'0.999"
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('D27').Value = "'2.70"
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('D28').Value = "'7.94"
$ws.Range('E28').Value = '  -2.00%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').Value = "'1.90"
$ws.Range('E30').Value = '  +1.27%  '
$ws.Range('D31').Value = "'25.64"
$ws.Range('E31').Value = '  -1.33%  '
$ws.Range('E32').Value = '  +4.01%  '
$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D33').Value = "'55.72"
$ws.Range('E33').Value = '  -2.38%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = "'5.49"
$ws.Range('E34').Value = '  +2.74%  '
$ws.Range('D35').Value = "'2.29"
$ws.Range('E35').Value = '  -3.20%  '
$ws.Range('D36').Value = "'5.90"
$ws.Range('E36').Value = '  -1.89%  '
$ws.Range('D37').Value = "'455.11"
$ws.Range('E37').Value = '  -6.32%  '
$ws.Range('D38').Value = '3.209.58'
$ws.Range('E38').Value = '  +2.69%  '
$ws.Range('D39').Value = "'0.0789"
$ws.Range('E39').Value = '  -0.69%  '
$ws.Range('D40').Value = "'0.0386"
$ws.Range('E40').Value = '  -1.97%  '
$ws.Range('D41').Value = "'0.118"
$ws.Range('E41').Value = '  +2.36%  '
$ws.Range('D42').Value = "'8.16"
$ws.Range('E42').Value = '  +0.86%  '
$ws.Range('D43').Value = "'27.85"
$ws.Range('E43').Value = '  +14.11%  '
$ws.Range('D44').Value = "'2.50"
$ws.Range('E44').Value = '  -3.93%  '
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('D46').Value = "'0.246"
$ws.Range('E46').Value = '  -2.48%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').Value = "'2.01"
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').Value = "'120.37"
$ws.Range('E48').Value = '  -1.06%  '
$ws.Range('D49').Value = "'0.108"
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('D50').Value = '0.0₃0499'
$ws.Range('E50').Value = '  -6.28%  '
$ws.Range('E51').Value = '  +7.25%  '
